# Update imputed values produced by re-running the RandomForest algorithm
# (commit: "Update Name of Algo")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -12.82309999999999
$ws.Range("C4").Value  = -14.1086
$ws.Range("D6").Value  = -7.551300000000003
$ws.Range("C7").Value  = -11.8331
$ws.Range("D7").Value  = -7.062999999999996
$ws.Range("C8").Value  = -11.91179999999999
$ws.Range("D8").Value  = -7.820299999999996
$ws.Range("B11").Value = 5.364300000000001
$ws.Range("E11").Value = 13.4273
$ws.Range("B12").Value = 5.5229
$ws.Range("C12").Value = -14.73390000000001
$ws.Range("C14").Value = -12.1103
$ws.Range("E14").Value = 14.0665
$ws.Range("B15").Value = 5.839100000000005
$ws.Range("D19").Value = -8.111099999999993
$ws.Range("E19").Value = 13.806
$ws.Range("D21").Value = -7.543000000000001
$ws.Range("E21").Value = 13.7893
$ws.Range("C22").Value = -11.79609999999999
$ws.Range("D24").Value = -7.522899999999995
$ws.Range("D25").Value = -7.574699999999996

$wb.Save()
